# Remove the "Diet Coke" row (row 6) from the vending items list.
# Excel shifts all subsequent rows up by one, drops the now-unused
# "Diet Coke" shared string, and shrinks the sheet's used range/dimension
# accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete() | Out-Null

# After the row delete Excel re-homes the view on A1 (the prior A8
# selection no longer makes sense once the sheet only spans A1:B8).
$ws.Range("A1").Select() | Out-Null
